$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '81.694.11'
Set-TextValue $ws 'D3' '3.197.62'
$ws.Range('E3').Value = '  +1.89%  '
Set-TextValue $ws 'D4' '1.00'
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue $ws 'D5' '210.06'
$ws.Range('E5').Value = '  +4.21%  '
Set-TextValue $ws 'D6' '637.58'
$ws.Range('E6').Value = '  +1.53%  '
$ws.Range('E7').Value = '  +28.46%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +3.91%  '
Set-TextValue $ws 'D10' '3.194.53'
$ws.Range('E10').Value = '  +1.84%  '
Set-TextValue $ws 'D11' '0.596'
$ws.Range('E11').Value = '  +11.60%  '
Set-TextValue $ws 'D12' '0.0000265'
$ws.Range('E12').Value = '  +18.70%  '
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('E14').Value = '  -1.04%  '
Set-TextValue $ws 'D15' '3.783.12'
$ws.Range('E15').Value = '  +1.78%  '
Set-TextValue $ws 'D16' '32.21'
$ws.Range('E16').Value = '  +5.70%  '
Set-TextValue $ws 'D17' '81.509.15'
$ws.Range('E17').Value = '  +5.63%  '
Set-TextValue $ws 'D18' '3.181.76'
$ws.Range('E18').Value = '  +1.33%  '
$ws.Range('E19').Value = '  +14.98%  '
Set-TextValue $ws 'D20' '14.42'
$ws.Range('E20').Value = '  +4.25%  '
Set-TextValue $ws 'D21' '9.33'
$ws.Range('E21').Value = '  +1.09%  '
Set-TextValue $ws 'D22' '443.46'
$ws.Range('E22').Value = '  +3.31%  '
Set-TextValue $ws 'D23' '5.25'
$ws.Range('E23').Value = '  +8.88%  '
Set-TextValue $ws 'D24' '7.12'
$ws.Range('E24').Value = '  +5.73%  '
$ws.Range('E25').Value = '  +9.96%  '
Set-TextValue $ws 'D26' '11.35'
$ws.Range('E26').Value = '  +6.08%  '
Set-TextValue $ws 'D27' '3.359.86'
$ws.Range('E27').Value = '  +1.81%  '
Set-TextValue $ws 'D28' '77.26'
$ws.Range('E28').Value = '  +2.39%  '
Set-TextValue $ws 'D29' '0.0000128'
$ws.Range('E29').Value = '  +12.44%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('E31').Value = '  +5.89%  '
Set-TextValue $ws 'D32' '0.999'
$ws.Range('E32').Value = '  +0.05%  '
Set-TextValue $ws 'D33' '575.30'
$ws.Range('E33').Value = '  +10.93%  '
$ws.Range('E34').Value = '  +2.22%  '
$ws.Range('E35').Value = '  +4.91%  '
Set-TextValue $ws 'D36' '0.154'
$ws.Range('E36').Value = '  +14.05%  '
Set-TextValue $ws 'D37' '0.142'
$ws.Range('E37').Value = '  +32.61%  '
Set-TextValue $ws 'D38' '23.41'
$ws.Range('E38').Value = '  +5.68%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws 'D39' '0.999'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws 'D40' '0.418'
$ws.Range('E40').Value = '  +6.55%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws 'D41' '2.08'
$ws.Range('E41').Value = '  +18.99%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 'D42' '3.12'
$ws.Range('E42').Value = '  +24.26%  '
Set-TextValue $ws 'D43' '5.99'
$ws.Range('E43').Value = '  +11.68%  '
$ws.Range('E44').Value = '  +3.72%  '
$ws.Range('E45').Value = '  -2.01%  '
Set-TextValue $ws 'D47' '189.51'
$ws.Range('E47').Value = '  -3.03%  '
Set-TextValue $ws 'D48' '45.39'
$ws.Range('E48').Value = '  +6.21%  '
$ws.Range('E49').Value = '  +6.10%  '
Set-TextValue $ws 'D50' '0.790'
$ws.Range('E50').Value = '  -0.92%  '
Set-TextValue $ws 'D51' '0.649'
$ws.Range('E51').Value = '  +5.41%  '
